$d = $word.ActiveDocument

# Locate the phrase being renamed ("constrained to " -> "categorized as ")
# include the trailing space so it is replaced exactly once.
$rng = $d.Content
$found = $rng.Find.Execute("constrained to ", $true, $false, $false, $false,
                            $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find target phrase 'constrained to '"
}

$start = $rng.Start
$end = $rng.End

# Remove the old wording.
$old = $d.Range($start, $end)
$old.Text = ""

# Insert the new wording as separate pieces -- "categorized" / " " / "as" / " " --
# each carrying the same (italic, non-bold, Courier New 9pt) run formatting
# that the surrounding sentence already uses.
$pieces = @("categorized", " ", "as", " ")
$pos = $start
$bounds = New-Object System.Collections.ArrayList
foreach ($piece in $pieces) {
    $insertionPoint = $d.Range($pos, $pos)
    $insertionPoint.InsertBefore($piece)
    $segEnd = $pos + $piece.Length
    [void]$bounds.Add(@($pos, $segEnd))
    $pos = $segEnd
}

foreach ($b in $bounds) {
    $seg = $d.Range($b[0], $b[1])
    $seg.Font.Name = "Courier New"
    $seg.Font.Italic = $true
    $seg.Font.Bold = $false
    $seg.Font.Size = 9
    $seg.Font.Underline = 0
}

Write-Output "Replaced 'constrained to' with 'categorized as'."
